$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 8).Value = 3.4   # H2
$ws.Cells.Item(2, 9).Value = 3.25   # I2
$ws.Cells.Item(2, 24).Value = 21   # X2
$ws.Cells.Item(2, 26).Value = 10   # Z2
$ws.Cells.Item(2, 27).Value = 7   # AA2
$ws.Cells.Item(2, 28).Value = 15   # AB2
$ws.Cells.Item(2, 36).Value = 41   # AJ2
$ws.Cells.Item(19, 14).Value = 1.44   # N19
$ws.Cells.Item(19, 15).Value = 2.7   # O19
$ws.Cells.Item(26, 30).Value = 101   # AD26
$ws.Cells.Item(35, 7).Value = 1.25   # G35
$ws.Cells.Item(35, 8).Value = 5   # H35
$ws.Cells.Item(35, 9).Value = 15   # I35
$ws.Cells.Item(35, 10).Value = 1.06   # J35
$ws.Cells.Item(35, 11).Value = 9.5   # K35
$ws.Cells.Item(35, 12).Value = 1.25   # L35
$ws.Cells.Item(35, 13).Value = 3.75   # M35
$ws.Cells.Item(35, 14).Value = 1.9   # N35
$ws.Cells.Item(35, 15).Value = 1.95   # O35
$ws.Cells.Item(35, 16).Value = 1.36   # P35
$ws.Cells.Item(35, 17).Value = 3   # Q35
$ws.Cells.Item(35, 18).Value = 2.5   # R35
$ws.Cells.Item(35, 19).Value = 1.5   # S35
$ws.Cells.Item(35, 20).Value = 5.5   # T35
$ws.Cells.Item(35, 21).Value = 5.5   # U35
$ws.Cells.Item(35, 23).Value = 7   # W35
$ws.Cells.Item(35, 24).Value = 13   # X35
$ws.Cells.Item(35, 26).Value = 9.5   # Z35
$ws.Cells.Item(35, 27).Value = 10   # AA35
$ws.Cells.Item(35, 28).Value = 29   # AB35
$ws.Cells.Item(35, 29).Value = 101   # AC35
$ws.Cells.Item(35, 31).Value = 23   # AE35
$ws.Cells.Item(35, 33).Value = 41   # AG35
$ws.Cells.Item(35, 34).Value = 201   # AH35
$ws.Cells.Item(37, 7).Value = 2.65   # G37
$ws.Cells.Item(37, 8).Value = 3.75   # H37
$ws.Cells.Item(37, 9).Value = 2.34   # I37
$ws.Cells.Item(37, 14).Value = 1.53   # N37
$ws.Cells.Item(37, 15).Value = 2.41   # O37
$ws.Cells.Item(37, 23).Value = 27   # W37
$ws.Cells.Item(37, 27).Value = 6.4   # AA37
$ws.Cells.Item(37, 28).Value = 9.4   # AB37
$ws.Cells.Item(37, 29).Value = 30   # AC37
$ws.Cells.Item(37, 33).Value = 8   # AG37
$ws.Cells.Item(37, 34).Value = 22   # AH37
$ws.Cells.Item(48, 7).Value = 8.5   # G48
$ws.Cells.Item(48, 8).Value = 4.9   # H48
$ws.Cells.Item(48, 9).Value = 1.27   # I48
$ws.Cells.Item(48, 12).Value = 1.19   # L48
$ws.Cells.Item(48, 13).Value = 4.2   # M48
$ws.Cells.Item(48, 14).Value = 1.57   # N48
$ws.Cells.Item(48, 15).Value = 2.1   # O48
$ws.Cells.Item(48, 18).Value = 2.09   # R48
$ws.Cells.Item(48, 19).Value = 1.66   # S48
$ws.Cells.Item(48, 20).Value = 18   # T48
$ws.Cells.Item(48, 21).Value = 45   # U48
$ws.Cells.Item(48, 22).Value = 22   # V48
$ws.Cells.Item(48, 23).Value = 175   # W48
$ws.Cells.Item(48, 24).Value = 80   # X48
$ws.Cells.Item(48, 25).Value = 70   # Y48
$ws.Cells.Item(48, 27).Value = 8.75   # AA48
$ws.Cells.Item(48, 28).Value = 19   # AB48
$ws.Cells.Item(48, 29).Value = 80   # AC48
$ws.Cells.Item(48, 31).Value = 6.2   # AE48
$ws.Cells.Item(48, 32).Value = 5.3   # AF48
$ws.Cells.Item(48, 33).Value = 7.5   # AG48
$ws.Cells.Item(48, 34).Value = 6.5   # AH48
$ws.Cells.Item(48, 36).Value = 23   # AJ48
$ws.Cells.Item(49, 7).Value = 3.7   # G49
$ws.Cells.Item(49, 8).Value = 3.35   # H49
$ws.Cells.Item(49, 9).Value = 1.85   # I49
$ws.Cells.Item(49, 12).Value = 1.29   # L49
$ws.Cells.Item(49, 13).Value = 3.3   # M49
$ws.Cells.Item(49, 14).Value = 1.8   # N49
$ws.Cells.Item(49, 15).Value = 1.8   # O49
$ws.Cells.Item(49, 16).Value = 1.37   # P49
$ws.Cells.Item(49, 17).Value = 2.5   # Q49
$ws.Cells.Item(49, 18).Value = 1.8   # R49
$ws.Cells.Item(49, 19).Value = 1.91   # S49
$ws.Cells.Item(49, 20).Value = 9.5   # T49
$ws.Cells.Item(49, 21).Value = 17   # U49
$ws.Cells.Item(49, 22).Value = 10.5   # V49
$ws.Cells.Item(49, 23).Value = 40   # W49
$ws.Cells.Item(49, 24).Value = 26   # X49
$ws.Cells.Item(49, 25).Value = 29   # Y49
$ws.Cells.Item(49, 26).Value = 10   # Z49
$ws.Cells.Item(49, 27).Value = 5.8   # AA49
$ws.Cells.Item(49, 28).Value = 11.75   # AB49
$ws.Cells.Item(49, 29).Value = 45   # AC49
$ws.Cells.Item(49, 31).Value = 6.2   # AE49
$ws.Cells.Item(49, 32).Value = 7.6   # AF49
$ws.Cells.Item(49, 33).Value = 7   # AG49
$ws.Cells.Item(49, 34).Value = 13   # AH49
$ws.Cells.Item(49, 36).Value = 20   # AJ49
$ws.Cells.Item(58, 7).Value = 1.42   # G58
$ws.Cells.Item(58, 9).Value = 7.5   # I58
$ws.Cells.Item(58, 10).Value = 1.06   # J58
$ws.Cells.Item(58, 11).Value = 10   # K58
$ws.Cells.Item(58, 22).Value = 8.5   # V58
$ws.Cells.Item(58, 23).Value = 9   # W58
$ws.Cells.Item(58, 25).Value = 29   # Y58
$ws.Cells.Item(58, 26).Value = 10   # Z58
$ws.Cells.Item(58, 28).Value = 21   # AB58
$ws.Cells.Item(58, 29).Value = 67   # AC58
$ws.Cells.Item(58, 34).Value = 81   # AH58
$ws.Cells.Item(58, 35).Value = 51   # AI58
$ws.Cells.Item(58, 36).Value = 51   # AJ58
$ws.Cells.Item(69, 10).Value = 1.05   # J69
$ws.Cells.Item(69, 11).Value = 11   # K69
$ws.Cells.Item(69, 14).Value = 2.03   # N69
$ws.Cells.Item(69, 15).Value = 1.83   # O69
$ws.Cells.Item(70, 7).Value = 2.55   # G70
$ws.Cells.Item(70, 8).Value = 3.6   # H70
$ws.Cells.Item(70, 9).Value = 2.55   # I70
$ws.Cells.Item(70, 10).Value = 1.04   # J70
$ws.Cells.Item(70, 11).Value = 12   # K70
$ws.Cells.Item(70, 12).Value = 1.25   # L70
$ws.Cells.Item(70, 13).Value = 3.75   # M70
$ws.Cells.Item(70, 14).Value = 1.88   # N70
$ws.Cells.Item(70, 15).Value = 1.98   # O70
$ws.Cells.Item(70, 16).Value = 1.36   # P70
$ws.Cells.Item(70, 17).Value = 3   # Q70
$ws.Cells.Item(70, 18).Value = 1.7   # R70
$ws.Cells.Item(70, 19).Value = 2.05   # S70
$ws.Cells.Item(70, 20).Value = 9   # T70
$ws.Cells.Item(70, 25).Value = 26   # Y70
$ws.Cells.Item(70, 26).Value = 12   # Z70
$ws.Cells.Item(70, 27).Value = 7   # AA70
$ws.Cells.Item(70, 28).Value = 13   # AB70
$ws.Cells.Item(70, 29).Value = 41   # AC70
$ws.Cells.Item(70, 30).Value = 201   # AD70
$ws.Cells.Item(70, 31).Value = 9   # AE70
$ws.Cells.Item(70, 36).Value = 26   # AJ70
$ws.Cells.Item(72, 7).Value = 2.42   # G72
$ws.Cells.Item(72, 8).Value = 3.25   # H72
$ws.Cells.Item(72, 9).Value = 2.75   # I72
$ws.Cells.Item(72, 12).Value = 1.22   # L72
$ws.Cells.Item(72, 13).Value = 4   # M72
$ws.Cells.Item(72, 14).Value = 1.71   # N72
$ws.Cells.Item(72, 15).Value = 2.06   # O72
$ws.Cells.Item(72, 18).Value = 1.58   # R72
$ws.Cells.Item(72, 19).Value = 2.3   # S72
$ws.Cells.Item(72, 20).Value = 8.2   # T72
$ws.Cells.Item(72, 21).Value = 11   # U72
$ws.Cells.Item(72, 22).Value = 7.6   # V72
$ws.Cells.Item(72, 23).Value = 20   # W72
$ws.Cells.Item(72, 24).Value = 15   # X72
$ws.Cells.Item(72, 25).Value = 21   # Y72
$ws.Cells.Item(72, 26).Value = 9.8   # Z72
$ws.Cells.Item(72, 27).Value = 5   # AA72
$ws.Cells.Item(72, 28).Value = 9   # AB72
$ws.Cells.Item(72, 29).Value = 35   # AC72
$ws.Cells.Item(72, 30).Value = 101   # AD72
$ws.Cells.Item(72, 31).Value = 8.8   # AE72
$ws.Cells.Item(72, 32).Value = 12   # AF72
$ws.Cells.Item(72, 33).Value = 8.2   # AG72
$ws.Cells.Item(72, 34).Value = 26   # AH72
$ws.Cells.Item(72, 35).Value = 17   # AI72
$ws.Cells.Item(72, 36).Value = 22   # AJ72
